$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "27.844.17"
Set-TextValue $ws.Range("E2") "  +3.26%  "
Set-TextValue $ws.Range("D3") "1.869.15"
Set-TextValue $ws.Range("E3") "  +2.80%  "
Set-TextValue $ws.Range("E4") "  +3.32%  "
Set-TextValue $ws.Range("D5") "325.30"
Set-TextValue $ws.Range("E5") "  +4.13%  "
Set-TextValue $ws.Range("D6") "1.038"
Set-TextValue $ws.Range("E6") "  +3.06%  "
Set-TextValue $ws.Range("D7") "0.4433"
Set-TextValue $ws.Range("E7") "  +3.06%  "
Set-TextValue $ws.Range("D8") "0.3806"
Set-TextValue $ws.Range("E8") "  +3.50%  "
Set-TextValue $ws.Range("D9") "0.07482"
Set-TextValue $ws.Range("E9") "  +3.21%  "
Set-TextValue $ws.Range("D10") "0.8871"
Set-TextValue $ws.Range("E10") "  +2.34%  "
Set-TextValue $ws.Range("D11") "21.83"
Set-TextValue $ws.Range("E11") "  +2.42%  "
Set-TextValue $ws.Range("D12") "1.885.83"
Set-TextValue $ws.Range("E12") "  -11.35%  "
Set-TextValue $ws.Range("D13") "5.569"
Set-TextValue $ws.Range("E13") "  +2.73%  "
Set-TextValue $ws.Range("D14") "6.772"
Set-TextValue $ws.Range("E14") "  +2.30%  "
Set-TextValue $ws.Range("D15") "0.07241"
Set-TextValue $ws.Range("E15") "  +3.68%  "
Set-TextValue $ws.Range("D16") "83.95"
Set-TextValue $ws.Range("E16") "  +3.37%  "
Set-TextValue $ws.Range("E17") "  +2.92%  "
Set-TextValue $ws.Range("D18") "0.000009195"
Set-TextValue $ws.Range("E18") "  +3.35%  "
Set-TextValue $ws.Range("D19") "1.038"
Set-TextValue $ws.Range("E19") "  +3.07%  "
Set-TextValue $ws.Range("D20") "15.60"
Set-TextValue $ws.Range("E20") "  +2.50%  "
Set-TextValue $ws.Range("D21") "27.851.19"
Set-TextValue $ws.Range("E21") "  +3.09%  "
Set-TextValue $ws.Range("D22") "5.331"
Set-TextValue $ws.Range("E22") "  +2.48%  "
Set-TextValue $ws.Range("E23") "  +3.41%  "
Set-TextValue $ws.Range("D24") "1.978"
Set-TextValue $ws.Range("E24") "  +4.62%  "
Set-TextValue $ws.Range("D25") "158.88"
Set-TextValue $ws.Range("E25") "  +3.05%  "
Set-TextValue $ws.Range("D26") "18.92"
Set-TextValue $ws.Range("E26") "  +2.79%  "
Set-TextValue $ws.Range("D27") "2.003"
Set-TextValue $ws.Range("E27") "  +5.20%  "
Set-TextValue $ws.Range("D28") "5.354"
Set-TextValue $ws.Range("E28") "  +2.41%  "
Set-TextValue $ws.Range("D29") "117.96"
Set-TextValue $ws.Range("E29") "  +2.67%  "
Set-TextValue $ws.Range("D30") "0.09105"
Set-TextValue $ws.Range("E30") "  +1.62%  "
Set-TextValue $ws.Range("D31") "0.7804"
Set-TextValue $ws.Range("E31") "  +4.13%  "
Set-TextValue $ws.Range("D32") "3.115"
Set-TextValue $ws.Range("E32") "  +10.71%  "
Set-TextValue $ws.Range("E33") "  +2.85%  "
Set-TextValue $ws.Range("D34") "4.586"
Set-TextValue $ws.Range("E34") "  +3.54%  "
Set-TextValue $ws.Range("D35") "1.040"
Set-TextValue $ws.Range("E35") "  +3.38%  "
Set-TextValue $ws.Range("D36") "1.159"
Set-TextValue $ws.Range("E36") "  +2.35%  "
Set-TextValue $ws.Range("D37") "0.01999"
Set-TextValue $ws.Range("E37") "  +3.84%  "
Set-TextValue $ws.Range("D38") "0.05362"
Set-TextValue $ws.Range("E38") "  +2.45%  "
Set-TextValue $ws.Range("D39") "2.867"
Set-TextValue $ws.Range("E39") "  +4.34%  "
Set-TextValue $ws.Range("D40") "0.5210"
Set-TextValue $ws.Range("E40") "  +1.98%  "
Set-TextValue $ws.Range("E41") "  +2.52%  "
Set-TextValue $ws.Range("D42") "6.923"
Set-TextValue $ws.Range("E42") "  +6.93%  "
Set-TextValue $ws.Range("D43") "8.695"
Set-TextValue $ws.Range("E43") "  +4.14%  "
Set-TextValue $ws.Range("D44") "110.30"
Set-TextValue $ws.Range("E44") "  +3.14%  "
Set-TextValue $ws.Range("D45") "10.77"
Set-TextValue $ws.Range("E45") "  +3.19%  "
Set-TextValue $ws.Range("E46") "  +5.27%  "
Set-TextValue $ws.Range("D47") "0.4721"
Set-TextValue $ws.Range("E47") "  +2.86%  "
Set-TextValue $ws.Range("D48") "0.06466"
Set-TextValue $ws.Range("E48") "  +3.83%  "
Set-TextValue $ws.Range("D49") "1.904"
Set-TextValue $ws.Range("E49") "  +4.01%  "
Set-TextValue $ws.Range("D50") "40.04"
Set-TextValue $ws.Range("E50") "  +3.96%  "
Set-TextValue $ws.Range("D51") "64.84"
Set-TextValue $ws.Range("E51") "  +2.63%  "
